# Update "想去人数" (want-to-go count) figures on both the "展览" sheet
# and the combined "全部类型" sheet, matching the upstream data refresh.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 162
$wsExpo.Range("F4").Value = 267
$wsExpo.Range("F5").Value = 4047
$wsExpo.Range("F7").Value = 449

# Sheet "全部类型" (All types) -- same events, but appear in different rows
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 162
$wsAll.Range("F4").Value = 267
$wsAll.Range("F5").Value = 4047
$wsAll.Range("F9").Value = 449
